$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear out the old data/styling entirely (B1:B3, A2:A3) so the sheet
# ends up containing only plain, unstyled text in A1:A2.
$ws.Range("A1:B3").Clear()

# Write the two file names as plain inline strings in column A.
$ws.Range("A1").Value = "Day 17 - The Ultimate TIP To Remembering Anything.MP3"
$ws.Range("A2").Value = "Day 17 - The Ultimate TIP To Remembering Anything.mp4"
